$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # detection_template_csv
$ws2 = $wb.Worksheets.Item(2)   # lookup

# ---------------------------------------------------------------------------
# Sheet "lookup": replace the old 1-2 column reference table with the new
# 4-column (term / website-label / ontological-definition / axiom) table.
# ---------------------------------------------------------------------------

# Clear everything first so no stale cells / widths survive.
$ws2.Cells.Clear()

$ws2.Range("A1").Value = "ClinEpi"
$ws2.Range("B1").Value = "ontological label"
$ws2.Range("C1").Value = "ontological definition"
$ws2.Range("D1").Value = "axiom"

$ws2.Range("A2").Value = "bacteriology"
$ws2.Range("B2").Value = "bacteriology"
$ws2.Range("C2").Value = "a bacteriological assay"
$ws2.Range("D2").Value = "bacteriological assay"

$ws2.Range("A3").Value = "PCR"
$ws2.Range("B3").Value = "PCR assay"
$ws2.Range("C3").Value = "an assay, of which a polymerase chain reaction is part,"
$ws2.Range("D3").Value = "assay and 'has part' some 'polymerase chain reaction'"

$ws2.Range("A4").Value = "TAC"
$ws2.Range("B4").Value = "TaqMan"
$ws2.Range("C4").Value = "a fluorogenic PCR assay"

# D4 is a formula that evaluates to a quoted label; the source cell was
# entered with a leading apostrophe (quote-prefix) format. Write the formula
# first, then build that format by typing a quote-prefixed value into a
# scratch cell and pasting its format (format only) onto D4 - pasting formats
# after setting the formula preserves the formula while picking up the
# quote-prefix xf.
$ws2.Range("D4").Formula = "=""'fluorogenic PCR assay'"""
$ws2.Range("ZZ1").Value = "'x"
$ws2.Range("ZZ1").Copy()
$ws2.Range("D4").PasteSpecial(-4122)
$ws2.Range("ZZ1").Clear()

# Column widths roughly matching the authored sheet.
$ws2.Columns.Item(2).ColumnWidth = 13.8
$ws2.Columns.Item(3).ColumnWidth = 46

# Persist a sort-state on the table (ascending by the term column), same as
# the author applied via Data > Sort.
$sort = $ws2.Sort()
$sort.SortFields().Clear()
$sort.SetRange($ws2.Range("A2:D6"))
$sort.SortFields().Add($ws2.Range("A2:A6"))
$sort.Header = 2
$sort.Apply()

$ws2.Range("E11").Select()

# ---------------------------------------------------------------------------
# Sheet "detection_template_csv": row 3 now refers to a "culture" assay
# (instead of "TAC"), the "by <assay>" formula uses the lookup table, and a
# scratch VLOOKUP demo was added in row 6.
# ---------------------------------------------------------------------------

$ws1.Range("C3").Value = "culture"

$ws1.Range("N3").Formula = "=IF(D3=""boolean"",""presence of"",IF(D3=""count"",""count of"",""data about""))&"" ""&H3&"" by ""&IF(ISNA(VLOOKUP(C3,lookup!A2:B4,2)=TRUE),C3,VLOOKUP(C3,lookup!A2:B4,2))"

$ws1.Range("M6").Formula = "=VLOOKUP(A3,lookup!A2:B4,2)"
$ws1.Range("J1").Copy()
$ws1.Range("M6").PasteSpecial(-4122)

$ws1.Range("M6").Select()
